# Hortaliza, Vega Central Mapocho de Santiago - Choclo
# Weekly data refresh: insert 3 new rows of data (latest date 44610 = 2022-02-18)
# above the former row 444, shifting the remaining rows (old 444:466) down to 447:469.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 444; this pushes the
# existing rows 444:466 down to 447:469 (matching the dimension growth
# from A1:R466 to A1:R469).
$ws.Rows("444:446").Insert()

# ---- New row 444 ----
$ws.Cells.Item(444, 1).Value2 = 9
$ws.Cells.Item(444, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(444, 3).Value2 = "Metropolitana"
$ws.Cells.Item(444, 4).Value2 = 44610
$ws.Cells.Item(444, 5).Value2 = 13
$ws.Cells.Item(444, 6).Value2 = 100112024
$ws.Cells.Item(444, 7).Value2 = "Choclo"
$ws.Cells.Item(444, 8).Value2 = "Choclero"
$ws.Cells.Item(444, 9).Value2 = "Primera"
$ws.Cells.Item(444, 10).Value2 = 7900
$ws.Cells.Item(444, 11).Value2 = 150
$ws.Cells.Item(444, 12).Value2 = 180
$ws.Cells.Item(444, 13).Value2 = 165
$ws.Cells.Item(444, 14).Value2 = "`$/unidad"
$ws.Cells.Item(444, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(444, 16).Value2 = 165
$ws.Cells.Item(444, 17).Value2 = 1
$ws.Cells.Item(444, 18).Value2 = "Hortaliza"

# ---- New row 445 ----
$ws.Cells.Item(445, 1).Value2 = 9
$ws.Cells.Item(445, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(445, 3).Value2 = "Metropolitana"
$ws.Cells.Item(445, 4).Value2 = 44610
$ws.Cells.Item(445, 5).Value2 = 13
$ws.Cells.Item(445, 6).Value2 = 100112024
$ws.Cells.Item(445, 7).Value2 = "Choclo"
$ws.Cells.Item(445, 8).Value2 = "Choclero"
$ws.Cells.Item(445, 9).Value2 = "Primera"
$ws.Cells.Item(445, 10).Value2 = 6100
$ws.Cells.Item(445, 11).Value2 = 150
$ws.Cells.Item(445, 12).Value2 = 180
$ws.Cells.Item(445, 13).Value2 = 165
$ws.Cells.Item(445, 14).Value2 = "`$/unidad"
$ws.Cells.Item(445, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(445, 16).Value2 = 165
$ws.Cells.Item(445, 17).Value2 = 1
$ws.Cells.Item(445, 18).Value2 = "Hortaliza"

# ---- New row 446 ----
$ws.Cells.Item(446, 1).Value2 = 9
$ws.Cells.Item(446, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(446, 3).Value2 = "Metropolitana"
$ws.Cells.Item(446, 4).Value2 = 44610
$ws.Cells.Item(446, 5).Value2 = 13
$ws.Cells.Item(446, 6).Value2 = 100112024
$ws.Cells.Item(446, 7).Value2 = "Choclo"
$ws.Cells.Item(446, 8).Value2 = "Dulce o Americano"
$ws.Cells.Item(446, 9).Value2 = "Primera"
$ws.Cells.Item(446, 10).Value2 = 4300
$ws.Cells.Item(446, 11).Value2 = 100
$ws.Cells.Item(446, 12).Value2 = 150
$ws.Cells.Item(446, 13).Value2 = 125
$ws.Cells.Item(446, 14).Value2 = "`$/unidad"
$ws.Cells.Item(446, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(446, 16).Value2 = 125
$ws.Cells.Item(446, 17).Value2 = 1
$ws.Cells.Item(446, 18).Value2 = "Hortaliza"

# Make sure the date cells keep/acquire the date-time display format used
# elsewhere in column D.
$ws.Range("D444:D446").NumberFormat = "YYYY-MM-DD HH:MM:SS"
